$wb = $excel.ActiveWorkbook

# Rename the "Exp" sheet to "LevelUp"
$ws = $wb.Worksheets.Item("Exp")
$ws.Name = "LevelUp"

# Make it the active/selected sheet (moves tabSelected + updates activeTab)
$ws.Activate()

# Update header cell A1 from "Id" to "Lev"
$ws.Range("A1").Value = "Lev"

# Touch PageSetup so a <pageSetup> element is emitted for this sheet
$ws.PageSetup.Orientation = 1
